# sigmas.xlsx - "changes made closer to finale"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Bump the multiplier (K column) from 1 to 2 for the first data
#    block (rows 2-13). Rows 14-19 stay at 1. H/I formulas recalc
#    automatically.
# ------------------------------------------------------------------
$ws.Range("K2:K13").Value = 2

# ------------------------------------------------------------------
# 2. Switch the computed "lower"/"upper" bound columns (H & I) to an
#    integer number format, and colour-code each analyte's rows with
#    a themed fill (different shade per analyte group).
# ------------------------------------------------------------------
$ws.Range("H2:I19").NumberFormat = "0"

# green leaf group (row 2) -> Accent 6, Darker 50%
$ws.Range("H2:I2").Interior.Color = 2381624

# darkgreen leaf group (row 5) -> Accent 6, Darker 25%
$ws.Range("H5:I5").Interior.Color = 3506772

# brown leaf group (row 8) -> Accent 6, Lighter 40%
$ws.Range("H8:I8").Interior.Color = 9555625

# lightgreen leaf group (row 11) -> Accent 6, Lighter 60%
$ws.Range("H11:I11").Interior.Color = 11919046

# dark dirt / dark brown dirt groups (rows 14 & 17) -> Accent 4, Darker 50%
$ws.Range("H14:I14").Interior.Color = 24704
$ws.Range("H17:I17").Interior.Color = 24704

# ------------------------------------------------------------------
# 3. Misc UI state: active selection moved, and the sheet now has
#    explicit page setup (paper size 9 = A4, portrait orientation).
# ------------------------------------------------------------------
$ws.Range("M8").Select()

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
